$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.844.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.081.46'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '551.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.14%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.073.99'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.66%  '
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.59'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.158'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.450'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.579.77'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.988.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.082.54'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '500.67'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.701'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.06%  '
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.95'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.11'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.23%  '
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.49'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '58.27'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +12.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '524.17'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.85'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.14'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0409'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.039.22'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.57%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.120'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.52%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0784'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.02'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.61'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.68%  '
$ws.Range('E44').Value = '  +4.47%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.106'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0498'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +69.68%  '
